$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update header C1: "CASH_KERAS" -> "HARGA_CASH_KERAS"
$ws.Range("C1").Value = "HARGA_CASH_KERAS"

# 2. Update column A (A2:A280): "MAR-2015" -> "FEB-2015"
# Leading apostrophe forces text entry (avoids Excel's date auto-detection)
$ws.Range("A2:A280").Value = "'FEB-2015"

# 3. Change selection to C5
$ws.Range("C5").Select()

# 4. Autofit column C to match new bestFit width
$ws.Columns("C").EntireColumn.AutoFit() | Out-Null
